$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'54.301.94"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -2.41%  '
$c = $ws.Range("D3")
$c.Value = "'2.288.43"
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.17%  '
$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.Value = "'494.26"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.66%  '
$c = $ws.Range("D6")
$c.Value = "'126.97"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -1.52%  '
$c = $ws.Range("D9")
$c.Value = "'2.289.65"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.38%  '
$c = $ws.Range("D10")
$c.Value = "'0.0939"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.12%  '
$ws.Range("E11").Value = '  +0.57%  '
$c = $ws.Range("D12")
$c.Value = "'0.320"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.53%  '
$c = $ws.Range("D13")
$c.Value = "'4.62"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.57%  '
$c = $ws.Range("D14")
$c.Value = "'2.695.98"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.18%  '
$c = $ws.Range("D15")
$c.Value = "'21.56"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.06%  '
$c = $ws.Range("D16")
$c.Value = "'54.139.58"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.69%  '
$c = $ws.Range("D17")
$c.Value = "'0.0000129"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.78%  '
$c = $ws.Range("D18")
$c.Value = "'2.295.52"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.15%  '
$c = $ws.Range("D19")
$c.Value = "'9.90"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.03%  '
$c = $ws.Range("D20")
$c.Value = "'4.04"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.66%  '
$c = $ws.Range("D21")
$c.Value = "'302.50"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.49%  '
$c = $ws.Range("D22")
$c.Value = "'6.43"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.57%  '
$c = $ws.Range("D23")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.09%  '
$c = $ws.Range("D24")
$c.Value = "'5.37"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.63%  '
$c = $ws.Range("D25")
$c.Value = "'63.72"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("E26").Value = '  +0.59%  '
$c = $ws.Range("D27")
$c.Value = "'0.374"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.12%  '
$c = $ws.Range("D28")
$c.Value = "'2.399.71"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.23%  '
$c = $ws.Range("D30")
$c.Value = "'7.12"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'165.40"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.56%  '
$ws.Range("E32").Value = '  -2.04%  '
$c = $ws.Range("D33")
$c.Value = "'0.0₃0681"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.64%  '
$c = $ws.Range("D34")
$c.Value = "'5.86"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.16%  '
$ws.Range("E35").Value = '  -0.01%  '
$c = $ws.Range("D36")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("E37").Value = '  +1.04%  '
$c = $ws.Range("D38")
$c.Value = "'17.58"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.18%  '
$c = $ws.Range("D39")
$c.Value = "'1.18"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.86%  '
$c = $ws.Range("D40")
$c.Value = "'0.864"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +5.73%  '
$c = $ws.Range("D41")
$c.Value = "'3.62"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.01%  '
$c = $ws.Range("D42")
$c.Value = "'35.43"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("E43").Value = '  +1.51%  '
$c = $ws.Range("D44")
$c.Value = "'1.39"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("E45").Value = '  +0.51%  '
$c = $ws.Range("D46")
$c.Value = "'126.09"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.35%  '
$c = $ws.Range("D47")
$c.Value = "'4.79"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.69%  '
$c = $ws.Range("D48")
$c.Value = "'0.0888"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '
$c = $ws.Range("D49")
$c.Value = "'0.545"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.37%  '
$c = $ws.Range("D50")
$c.Value = "'236.71"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.19%  '
$c = $ws.Range("D51")
$c.Value = "'0.0480"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.49%  '
